# "Add Gen Dupe List" - insert a new "ListOfDupes" worksheet between
# "Report" and "Formats" that lists team members (col B) against a
# second list containing letters + some repeated names (col C), then
# uses FILTER/COUNTIF array formulas to surface the values that occur
# in both lists (i.e. the "dupes").

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "Report" so the tab order becomes
# Report, ListOfDupes, Formats, Lists.
$reportSheet = $wb.Worksheets.Item("Report")
$newSheet = $wb.Worksheets.Add($null, $reportSheet)
$newSheet.Name = "ListOfDupes"

# Standard memo header block (reuses the same styles as the other
# sheets' "FROM:/SUBJECT:/DATE:" headers).
$newSheet.Range("A1:A3").Style = "Intro_Hd"
$newSheet.Range("B1:C3").Style = "Intro_Value"

$newSheet.Range("A1").Value = "FROM:"
$newSheet.Range("B1").Value = "Mark Biegert"
$newSheet.Range("A2").Value = "SUBJECT:"
$newSheet.Range("B2").Value = "Dupe Methods"
$newSheet.Range("A3").Value = "DATE:"
$newSheet.Range("B3").FormulaArray = "=HD_Date"

# Column B: list of names (B8:B17)
$names = @("Mark","Craig","Tim","Pat","Katee","Biff","Zip","Dale","Miki","Thia")
for ($i = 0; $i -lt $names.Length; $i++) {
  $row = 8 + $i
  $newSheet.Cells.Item($row, 2).Value = $names[$i]
}

# Column C: letters A-F followed by a few repeated names (C8:C18)
$letters = @("A","B","C","D","E","F","Mark","Craig","Tim","Pat","Katee")
for ($i = 0; $i -lt $letters.Length; $i++) {
  $row = 8 + $i
  $newSheet.Cells.Item($row, 3).Value = $letters[$i]
}

# Dupe-finder array formulas.
$newSheet.Range("D8:D12").FormulaArray = "=FILTER(B8:B17,COUNTIF(C8:C18,B8:B17))"
$newSheet.Range("E8:E12").FormulaArray = "=FILTER(C8:C18,COUNTIF(B8:B17,C8:C18))"

# Match the selection/active-cell the author left on the new sheet.
$newSheet.Range("G9").Select()

# Best-effort: restore window placement recorded in the original
# session (not all hosts persist window geometry into the saved file).
$win = $excel.ActiveWindow
$win.Left = 9765
$win.Top = 1020
$win.Width = 14160
$win.Height = 15240
